$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2017-02-15 06:04:06"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2017-02-15 06:03:48"
$wsZhCn.Range("L2").Value = "2017-02-15 06:04:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2017-02-15 06:04:06"
$wsDeDe.Range("L2").Value = "2017-02-15 06:05:09"
